$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.618.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.560"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.70"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.74%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.15%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.079.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.852.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.65%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.33"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.644.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.90%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.63"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.90%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.26%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.79%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.02%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.07%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.418.15"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.675"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.48%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "86.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.34%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.959"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.26%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0526"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.71%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.980.04"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.97"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.56%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.00%  "
